$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.571.94"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.974.59"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'244.40"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'0.621"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").Value = "'60.03"
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("D10").Value = "'0.0790"
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("E12").Value = "  +3.47%  "
$ws.Range("D13").Value = "'0.842"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "2.263.92"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "'21.66"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "'5.41"
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("D17").Value = "1.980.57"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "36.555.68"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'69.82"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'229.60"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.09"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  +7.85%  "
$ws.Range("D27").Value = "'9.17"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "'1.35"
$ws.Range("E30").Value = "  +20.10%  "
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").Value = "'4.50"
$ws.Range("E34").Value = "  +5.97%  "
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.77"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'3.29"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("E39").Value = "  -11.95%  "
$ws.Range("D40").Value = "'0.0972"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").Value = "'15.93"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "1.366.63"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "'89.15"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").Value = "'46.24"
$ws.Range("E50").Value = "  +6.62%  "
$ws.Range("D51").Value = "2.157.94"
$ws.Range("E51").Value = "  +0.50%  "
